$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990868838/1749474758173/gFRHQATQ91-camera-video-segment-driver-1749474760673.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSASICHJTH%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T210003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEO3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIEOnqJwrJ3K3rIds%2B1hSvBAx0P6uCW8fmEDBeK8bN7shAiEA5Rf%2F0%2FQ91vmgWtPXTL94iFjhFD5cYu5fp3NNjQWTGzUq5gMIxv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDMR511Fy6RuuxJtPkyq6A5WnZVmtytURnaggzCAVDnemobLWSBthLiMVsQf%2F6dlwYxveoQOi1Z5%2BhWyyTJkW7cDxS4%2B2t7vSUl02714OWIorV6OTi6ZVVycxx5SBXeeMvy2XrNqII4Addmv6ZR21HwEmH2RsOfPLFu8g6EZm3ZGuUSVoprGEwqUjVoqrQ0wyaq%2Fmo45%2BZi5r7FGF61x%2BvWOH36fu5mTu4fpgBMmYnIBhbNQF2buSPdSVrWMUT%2BTovdDEEzWJRxfFdNnfOmScMIVzo9R%2Fn6SYo5gU33MFOCGPnV9W3QDds1u5iqUQUBodW5zULGJ4TuMUgEKtbOQqcESaddaHUV2A1L%2FPhd9kv888JFY4%2B6kEdVJWEQok17QJctwA0fb7BeuOhxgaWHsKuGa4sXmERXiLt0r9dkbe%2FfXsWmPI7BiAgCZZS8nH5uv9vW1YcmdJRxWwxGYqGRCwsG84jy5v95HSVP%2FGLQV2QiJa%2BZiGcMS3Y%2BHH5tuMruHX0M8MQ%2BM3z7Y4HW3n8LGOK83m4e%2FkauYQscUFjGnMRDnmxEaHAt0G0eXNBCJWvLVwEJ4DcLDMtoIwwL%2BZ5YRlW%2FWqksiSk%2BwQeQgwkquiwgY6pQFdWfclZnD%2F4Zwz9Ou97DMXsjYfs2WSAjjxHfvx2%2BRvbrm2VogWOT01veYp2C2omL2QdtswM7k7CgMwYAb1dFJaC8V1VVGXyxrCZ8Ror%2F3jfrqbISx7IFBcyl49kBhhArPxO1kMNOhl8PJ82kDIXxdz%2BBw6DMGSOl6xgd3iD%2F5AcyL7mLzwS3PAxZW64iuavWu166q%2FJu5XSqfofpwM4%2F58oHj9rzs%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2005%3A00%3A03%20GMT&X-Amz-Signature=cb8d860983e4d29620331e72cd0af7932f8e98a353191d89ed63f38a25ab2290'

$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990868838/1749474758173/7Z3ova6Wtr-camera-video-segment-1749474760673.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSASICHJTH%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T210003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEO3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIEOnqJwrJ3K3rIds%2B1hSvBAx0P6uCW8fmEDBeK8bN7shAiEA5Rf%2F0%2FQ91vmgWtPXTL94iFjhFD5cYu5fp3NNjQWTGzUq5gMIxv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDMR511Fy6RuuxJtPkyq6A5WnZVmtytURnaggzCAVDnemobLWSBthLiMVsQf%2F6dlwYxveoQOi1Z5%2BhWyyTJkW7cDxS4%2B2t7vSUl02714OWIorV6OTi6ZVVycxx5SBXeeMvy2XrNqII4Addmv6ZR21HwEmH2RsOfPLFu8g6EZm3ZGuUSVoprGEwqUjVoqrQ0wyaq%2Fmo45%2BZi5r7FGF61x%2BvWOH36fu5mTu4fpgBMmYnIBhbNQF2buSPdSVrWMUT%2BTovdDEEzWJRxfFdNnfOmScMIVzo9R%2Fn6SYo5gU33MFOCGPnV9W3QDds1u5iqUQUBodW5zULGJ4TuMUgEKtbOQqcESaddaHUV2A1L%2FPhd9kv888JFY4%2B6kEdVJWEQok17QJctwA0fb7BeuOhxgaWHsKuGa4sXmERXiLt0r9dkbe%2FfXsWmPI7BiAgCZZS8nH5uv9vW1YcmdJRxWwxGYqGRCwsG84jy5v95HSVP%2FGLQV2QiJa%2BZiGcMS3Y%2BHH5tuMruHX0M8MQ%2BM3z7Y4HW3n8LGOK83m4e%2FkauYQscUFjGnMRDnmxEaHAt0G0eXNBCJWvLVwEJ4DcLDMtoIwwL%2BZ5YRlW%2FWqksiSk%2BwQeQgwkquiwgY6pQFdWfclZnD%2F4Zwz9Ou97DMXsjYfs2WSAjjxHfvx2%2BRvbrm2VogWOT01veYp2C2omL2QdtswM7k7CgMwYAb1dFJaC8V1VVGXyxrCZ8Ror%2F3jfrqbISx7IFBcyl49kBhhArPxO1kMNOhl8PJ82kDIXxdz%2BBw6DMGSOl6xgd3iD%2F5AcyL7mLzwS3PAxZW64iuavWu166q%2FJu5XSqfofpwM4%2F58oHj9rzs%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2005%3A00%3A03%20GMT&X-Amz-Signature=05d6a2ac7ccbda2752f98de16571e34ae5e577717dbc83de09362edd262c161b'

$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990868838/1749474688613/JLWgs3viB6-camera-video-segment-driver-1749474691113.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSASICHJTH%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T210003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEO3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIEOnqJwrJ3K3rIds%2B1hSvBAx0P6uCW8fmEDBeK8bN7shAiEA5Rf%2F0%2FQ91vmgWtPXTL94iFjhFD5cYu5fp3NNjQWTGzUq5gMIxv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDMR511Fy6RuuxJtPkyq6A5WnZVmtytURnaggzCAVDnemobLWSBthLiMVsQf%2F6dlwYxveoQOi1Z5%2BhWyyTJkW7cDxS4%2B2t7vSUl02714OWIorV6OTi6ZVVycxx5SBXeeMvy2XrNqII4Addmv6ZR21HwEmH2RsOfPLFu8g6EZm3ZGuUSVoprGEwqUjVoqrQ0wyaq%2Fmo45%2BZi5r7FGF61x%2BvWOH36fu5mTu4fpgBMmYnIBhbNQF2buSPdSVrWMUT%2BTovdDEEzWJRxfFdNnfOmScMIVzo9R%2Fn6SYo5gU33MFOCGPnV9W3QDds1u5iqUQUBodW5zULGJ4TuMUgEKtbOQqcESaddaHUV2A1L%2FPhd9kv888JFY4%2B6kEdVJWEQok17QJctwA0fb7BeuOhxgaWHsKuGa4sXmERXiLt0r9dkbe%2FfXsWmPI7BiAgCZZS8nH5uv9vW1YcmdJRxWwxGYqGRCwsG84jy5v95HSVP%2FGLQV2QiJa%2BZiGcMS3Y%2BHH5tuMruHX0M8MQ%2BM3z7Y4HW3n8LGOK83m4e%2FkauYQscUFjGnMRDnmxEaHAt0G0eXNBCJWvLVwEJ4DcLDMtoIwwL%2BZ5YRlW%2FWqksiSk%2BwQeQgwkquiwgY6pQFdWfclZnD%2F4Zwz9Ou97DMXsjYfs2WSAjjxHfvx2%2BRvbrm2VogWOT01veYp2C2omL2QdtswM7k7CgMwYAb1dFJaC8V1VVGXyxrCZ8Ror%2F3jfrqbISx7IFBcyl49kBhhArPxO1kMNOhl8PJ82kDIXxdz%2BBw6DMGSOl6xgd3iD%2F5AcyL7mLzwS3PAxZW64iuavWu166q%2FJu5XSqfofpwM4%2F58oHj9rzs%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2005%3A00%3A03%20GMT&X-Amz-Signature=85e2eec8a06a4faea431dfd4acba4d28f35ce9b65ac96f6558818e17f62b1b02'

$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990868838/1749474688613/LOHXs2uAgt-camera-video-segment-1749474691113.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSASICHJTH%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T210003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEO3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIEOnqJwrJ3K3rIds%2B1hSvBAx0P6uCW8fmEDBeK8bN7shAiEA5Rf%2F0%2FQ91vmgWtPXTL94iFjhFD5cYu5fp3NNjQWTGzUq5gMIxv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDMR511Fy6RuuxJtPkyq6A5WnZVmtytURnaggzCAVDnemobLWSBthLiMVsQf%2F6dlwYxveoQOi1Z5%2BhWyyTJkW7cDxS4%2B2t7vSUl02714OWIorV6OTi6ZVVycxx5SBXeeMvy2XrNqII4Addmv6ZR21HwEmH2RsOfPLFu8g6EZm3ZGuUSVoprGEwqUjVoqrQ0wyaq%2Fmo45%2BZi5r7FGF61x%2BvWOH36fu5mTu4fpgBMmYnIBhbNQF2buSPdSVrWMUT%2BTovdDEEzWJRxfFdNnfOmScMIVzo9R%2Fn6SYo5gU33MFOCGPnV9W3QDds1u5iqUQUBodW5zULGJ4TuMUgEKtbOQqcESaddaHUV2A1L%2FPhd9kv888JFY4%2B6kEdVJWEQok17QJctwA0fb7BeuOhxgaWHsKuGa4sXmERXiLt0r9dkbe%2FfXsWmPI7BiAgCZZS8nH5uv9vW1YcmdJRxWwxGYqGRCwsG84jy5v95HSVP%2FGLQV2QiJa%2BZiGcMS3Y%2BHH5tuMruHX0M8MQ%2BM3z7Y4HW3n8LGOK83m4e%2FkauYQscUFjGnMRDnmxEaHAt0G0eXNBCJWvLVwEJ4DcLDMtoIwwL%2BZ5YRlW%2FWqksiSk%2BwQeQgwkquiwgY6pQFdWfclZnD%2F4Zwz9Ou97DMXsjYfs2WSAjjxHfvx2%2BRvbrm2VogWOT01veYp2C2omL2QdtswM7k7CgMwYAb1dFJaC8V1VVGXyxrCZ8Ror%2F3jfrqbISx7IFBcyl49kBhhArPxO1kMNOhl8PJ82kDIXxdz%2BBw6DMGSOl6xgd3iD%2F5AcyL7mLzwS3PAxZW64iuavWu166q%2FJu5XSqfofpwM4%2F58oHj9rzs%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2005%3A00%3A03%20GMT&X-Amz-Signature=d497e5880b71cc59832d5b5404efdbe8bd869f2a15b485353b44e74e925c11ae'

